$wb = $excel.ActiveWorkbook

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1420.2858
$ws.Range("J88").Value = 1420.2858
$ws.Range("L88").Value = 1420.2858
$ws.Range("N88").Value = -2232.2858

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1420.2858
$ws.Range("J91").Value = 1420.2858
$ws.Range("L91").Value = 1420.2858
$ws.Range("N91").Value = -4228.2858

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2365.3403
$ws.Range("I137").Value = 1846.6316
$ws.Range("J137").Value = 4555.4443
$ws.Range("K137").Value = 5539.8948
$ws.Range("L137").Value = 13666.3329
$ws.Range("M137").Value = -2989.8948
$ws.Range("N137").Value = -18766.3329

# ARM row 18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30644

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9463.375
$ws.Range("I61").Value = 6157.7856
$ws.Range("J61").Value = 17176.416
$ws.Range("K61").Value = 6157.7856
$ws.Range("L61").Value = 17176.416
$ws.Range("M61").Value = -5945.7856
$ws.Range("N61").Value = -17600.416

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 85417.33
$ws.Range("I74").Value = 96287.484
$ws.Range("J74").Value = 18384.666
$ws.Range("K74").Value = 96287.484
$ws.Range("L74").Value = 18384.666
$ws.Range("M74").Value = -95413.484
$ws.Range("N74").Value = -20132.666

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 85417.33
$ws.Range("I77").Value = 96287.484
$ws.Range("J77").Value = 18384.666
$ws.Range("K77").Value = 481437.42
$ws.Range("L77").Value = 91923.33
$ws.Range("M77").Value = -477069.42
$ws.Range("N77").Value = -100659.33

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4674.237
$ws.Range("I132").Value = 4842.2256
$ws.Range("J132").Value = 3930.2856
$ws.Range("K132").Value = 14526.6768
$ws.Range("L132").Value = 11790.8568
$ws.Range("M132").Value = -11996.6768
$ws.Range("N132").Value = -16850.8568

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 9463.375
$ws.Range("I136").Value = 6157.7856
$ws.Range("J136").Value = 17176.416
$ws.Range("K136").Value = 18473.3568
$ws.Range("L136").Value = 51529.24800000001
$ws.Range("M136").Value = -15923.3568
$ws.Range("N136").Value = -56629.24800000001

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 36961.414
$ws.Range("I134").Value = 2535.7727
$ws.Range("J134").Value = 145156.28
$ws.Range("K134").Value = 7607.3181
$ws.Range("L134").Value = 435468.84
$ws.Range("M134").Value = -5072.3181
$ws.Range("N134").Value = -440538.84

# CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 102.666664
$ws.Range("I3").Value = 98
$ws.Range("K3").Value = 98
$ws.Range("M3").Value = 15

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3078.8064
$ws.Range("I31").Value = 2729.9092
$ws.Range("K31").Value = 2729.9092
$ws.Range("M31").Value = -2434.9092

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3078.8064
$ws.Range("I34").Value = 2729.9092
$ws.Range("K34").Value = 2729.9092
$ws.Range("M34").Value = -2527.9092

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2599903.5
$ws.Range("I58").Value = 5683275.5
$ws.Range("J58").Value = 3379.6843
$ws.Range("K58").Value = 5683275.5
$ws.Range("L58").Value = 3379.6843
$ws.Range("M58").Value = -5683072.5
$ws.Range("N58").Value = -3785.6843

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2508.721
$ws.Range("I132").Value = 2325.7942
$ws.Range("J132").Value = 3199.7778
$ws.Range("K132").Value = 6977.382599999999
$ws.Range("L132").Value = 9599.3334
$ws.Range("M132").Value = -4447.382599999999
$ws.Range("N132").Value = -14659.3334

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3088.8696
$ws.Range("I134").Value = 3079.8235
$ws.Range("J134").Value = 3114.5
$ws.Range("K134").Value = 9239.470499999999
$ws.Range("L134").Value = 9343.5
$ws.Range("M134").Value = -6704.470499999999
$ws.Range("N134").Value = -14413.5

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2599903.5
$ws.Range("I136").Value = 5683275.5
$ws.Range("J136").Value = 3379.6843
$ws.Range("K136").Value = 17049826.5
$ws.Range("L136").Value = 10139.0529
$ws.Range("M136").Value = -17047276.5
$ws.Range("N136").Value = -15239.0529

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 19706.889
$ws.Range("J141").Value = 15235.6
$ws.Range("L141").Value = 15235.6
$ws.Range("N141").Value = -25595.6

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1759
$ws.Range("I117").Value = 1118
$ws.Range("J117").Value = 2400
$ws.Range("K117").Value = 3354
$ws.Range("L117").Value = 7200
$ws.Range("M117").Value = 88
$ws.Range("N117").Value = -14084

# CUL row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 8333.333000000001
$ws.Range("I120").Value = 14000
$ws.Range("J120").Value = 5500
$ws.Range("K120").Value = 42000
$ws.Range("L120").Value = 16500
$ws.Range("M120").Value = -37162
$ws.Range("N120").Value = -26176

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20806.125
$ws.Range("J131").Value = 29021.588
$ws.Range("L131").Value = 87064.764
$ws.Range("N131").Value = -97144.764

# GSM row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 116780280
$ws.Range("J3").Value = 4249.75
$ws.Range("L3").Value = 4249.75
$ws.Range("N3").Value = -4481.75

# GSM row 7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10000033
$ws.Range("I7").Value = 10000050
$ws.Range("K7").Value = 10000050
$ws.Range("M7").Value = -9999938

# GSM row 8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 10000033
$ws.Range("I8").Value = 10000050
$ws.Range("K8").Value = 10000050
$ws.Range("M8").Value = -9999911

# GSM row 9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 8450
$ws.Range("I9").Value = 2900
$ws.Range("J9").Value = 14000
$ws.Range("K9").Value = 2900
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = -2730
$ws.Range("N9").Value = -14340

# GSM row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J10").Value = 9999
$ws.Range("L10").Value = 9999
$ws.Range("N10").Value = -10337

# GSM row 13
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 2167.9
$ws.Range("I13").Value = 1462.375
$ws.Range("J13").Value = 4990
$ws.Range("K13").Value = 1462.375
$ws.Range("L13").Value = 4990
$ws.Range("M13").Value = -1323.375
$ws.Range("N13").Value = -5268

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 96456790
$ws.Range("I14").Value = 151572100
$ws.Range("K14").Value = 151572100
$ws.Range("M14").Value = -151571932

# GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3999
$ws.Range("J17").Value = 3999
$ws.Range("L17").Value = 3999
$ws.Range("N17").Value = -4335

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 6015120.5
$ws.Range("J18").Value = 18899.5
$ws.Range("L18").Value = 18899.5
$ws.Range("N18").Value = -19485.5

# GSM row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4980
$ws.Range("I19").Value = 4980
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 4980
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -4692
$ws.Range("N19").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4657.5713
$ws.Range("I80").Value = 2540
$ws.Range("J80").Value = 5834
$ws.Range("K80").Value = 2540
$ws.Range("L80").Value = 5834
$ws.Range("M80").Value = -1542
$ws.Range("N80").Value = -7830

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4657.5713
$ws.Range("I83").Value = 2540
$ws.Range("J83").Value = 5834
$ws.Range("K83").Value = 12700
$ws.Range("L83").Value = 29170
$ws.Range("M83").Value = -7708
$ws.Range("N83").Value = -39154

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6494.5
$ws.Range("I132").Value = 2649.7058
$ws.Range("J132").Value = 13756.889
$ws.Range("K132").Value = 7949.117400000001
$ws.Range("L132").Value = 41270.667
$ws.Range("M132").Value = -5419.117400000001
$ws.Range("N132").Value = -46330.667

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4648.4287
$ws.Range("I7").Value = 4600.5713
$ws.Range("J7").Value = 4744.143
$ws.Range("K7").Value = 4600.5713
$ws.Range("L7").Value = 4744.143
$ws.Range("M7").Value = -4488.5713
$ws.Range("N7").Value = -4968.143

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6559.551
$ws.Range("I122").Value = 6036.2256
$ws.Range("J122").Value = 7460.8335
$ws.Range("K122").Value = 18108.6768
$ws.Range("L122").Value = 22382.5005
$ws.Range("M122").Value = -15658.6768
$ws.Range("N122").Value = -27282.5005

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4648.4287
$ws.Range("I126").Value = 4600.5713
$ws.Range("J126").Value = 4744.143
$ws.Range("K126").Value = 13801.7139
$ws.Range("L126").Value = 14232.429
$ws.Range("M126").Value = -11331.7139
$ws.Range("N126").Value = -19172.429

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4459
$ws.Range("I132").Value = 4793.967
$ws.Range("J132").Value = 3789.0667
$ws.Range("K132").Value = 14381.901
$ws.Range("L132").Value = 11367.2001
$ws.Range("M132").Value = -11851.901
$ws.Range("N132").Value = -16427.2001

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4138.1963
$ws.Range("I136").Value = 2417.3333
$ws.Range("K136").Value = 7251.999899999999
$ws.Range("M136").Value = -4701.999899999999

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2681.818
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2681.818
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2681.818
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2907.818

# WVR row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20005000
$ws.Range("J5").Value = 20005000
$ws.Range("L5").Value = 20005000
$ws.Range("N5").Value = -20005224

# WVR row 6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3622.2
$ws.Range("I6").Value = 368.33334
$ws.Range("J6").Value = 8503
$ws.Range("K6").Value = 368.33334
$ws.Range("L6").Value = 8503
$ws.Range("M6").Value = -253.33334
$ws.Range("N6").Value = -8733

# WVR row 8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# WVR row 11
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 3750
$ws.Range("I11").Value = 2500
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -2358
$ws.Range("N11").Value = -5284

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1722.9546
$ws.Range("I132").Value = 1131.6786
$ws.Range("J132").Value = 2757.6875
$ws.Range("K132").Value = 3395.0358
$ws.Range("L132").Value = 8273.0625
$ws.Range("M132").Value = -865.0357999999997
$ws.Range("N132").Value = -13333.0625

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8539.843999999999
$ws.Range("I136").Value = 8467.111000000001
$ws.Range("J136").Value = 8568.305
$ws.Range("K136").Value = 25401.333
$ws.Range("L136").Value = 25704.915
$ws.Range("M136").Value = -22851.333
$ws.Range("N136").Value = -30804.915
